$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test-result")

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 10).Value = "2025-06-07 13:13:37"
}
